$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.297.50'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.506.89'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').Value = '2.527.47'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0981'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.16'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('E13').Value = '  -2.89%  '
$ws.Range('D14').Value = '2.973.05'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '58.246.90'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '2.524.55'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.89'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  -1.83%  '
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.994'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.21'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('E37').Value = '  -6.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.93'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.49'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.773'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '277.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('E44').Value = '  -4.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '129.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.599'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0919'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0500'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.79%  '
